# Slide 17 ("Optimization in CPRL") - Content Placeholder 7:
#   - Reword the intro sentence about the assembler/peephole optimizations.
#   - Swap the order of the "constant folding" and "branch reduction" bullets,
#     and split "constant folding" into two runs ("constant " + "folding").

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(17)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame.TextRange

$quoteOpen  = [char]0x201C
$quoteClose = [char]0x201D

# "The assembler for CPRL performs ..." -> "The CVM assembler performs ..."
$tf.Paragraphs(4).Runs(1).Text = "The CVM assembler performs the following optimizations using a " + $quoteOpen + "peephole" + $quoteClose + " approach:"

# Bullet that used to read "constant folding" now reads "branch reduction ..."
$tf.Paragraphs(5).Runs(1).Text = "branch reduction (as illustrated in previous slide)"

# Bullet that used to read "branch reduction ..." now reads "constant folding",
# split across two runs: "constant " and "folding".
$tf.Paragraphs(6).Runs(1).Text = "constant folding"
$para6 = $tf.Paragraphs(6)
$secondRun = $para6.Characters(10, 7)
$secondRun.Text = "folding"
